# Generate Report for Handback
# Refreshes the Correspond Handoff/Handback datetimes for the
# "bfc56466-f424-4c71-a9e5-4b645e843490" file row on each language sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 00:51:31"
$wsZhCn.Range("H2").Value = "2016-03-24 00:51:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 00:51:35"
$wsDeDe.Range("H2").Value = "2016-03-24 00:52:03"
